$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format so values like "1.000" or
# "29.401.26" are preserved exactly instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$data = @(
    ,@('Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '29.401.26', '  -1.07%  ')
    ,@('Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.896.22', '  -1.43%  ')
    ,@('TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.001', '  +0.02%  ')
    ,@('BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '323.92', '  -3.46%  ')
    ,@('USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.000', '  +0.04%  ')
    ,@('XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.4765', '  +1.83%  ')
    ,@('Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.4055', '  -2.13%  ')
    ,@('OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '47.48', '  -1.56%  ')
    ,@('Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.08014', '  -0.58%  ')
    ,@('Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '1.001', '  -1.59%  ')
    ,@('Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '23.38', '  +4.33%  ')
    ,@('Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.922', '  -1.73%  ')
    ,@('WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.828.16', '  -4.35%  ')
    ,@('Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '7.052', '  -2.23%  ')
    ,@('Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '89.48', '  -0.66%  ')
    ,@('BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.001', '  +0.04%  ')
    ,@('TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.06684', '  +1.45%  ')
    ,@('ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.00001029', '  -0.81%  ')
    ,@('Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '17.57', '  -1.70%  ')
    ,@('Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.000', '  +0.09%  ')
    ,@('WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '29.417.24', '  -0.95%  ')
    ,@('Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '5.517', '  -0.90%  ')
    ,@('Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '11.68', '  +0.34%  ')
    ,@('Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '2.157', '  -1.97%  ')
    ,@('WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.136.21', '  -0.76%  ')
    ,@('Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '154.47', '  -1.59%  ')
    ,@('EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '19.76', '  -1.01%  ')
    ,@('InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '6.029', '  +5.30%  ')
    ,@('LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '2.084', '  -3.13%  ')
    ,@('BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '117.76', '  -0.14%  ')
    ,@('ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '1.018', '  -4.18%  ')
    ,@('Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.09458', '  -0.15%  ')
    ,@('ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.386', '  -3.64%  ')
    ,@('HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '3.527', '  +0.13%  ')
    ,@('Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '5.351', '  -1.60%  ')
    ,@('VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.02248', '  -1.11%  ')
    ,@('Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.06037', '  -2.01%  ')
    ,@('TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '1.167', '  -1.43%  ')
    ,@('TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.5855', '  -1.15%  ')
    ,@('FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '7.822', '  -7.55%  ')
    ,@('Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1836', '  -0.73%  ')
    ,@('Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '10.10', '  -1.87%  ')
    ,@('RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '2.423', '  +3.19%  ')
    ,@('WEMIXToken', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix', '1.286', '  +2.31%  ')
    ,@('Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.07711', '  +2.49%  ')
    ,@('EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '12.18', '  -0.26%  ')
    ,@('Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.5489', '  -2.09%  ')
    ,@('NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '1.915', '  -1.51%  ')
    ,@('Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '112.94', '  -0.06%  ')
    ,@('WOONetwork', 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo', '0.2979', '  -0.47%  ')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value = $data[$i][0]
    $ws.Cells.Item($r, 3).Value = $data[$i][1]
    $ws.Cells.Item($r, 4).Value = $data[$i][2]
    $ws.Cells.Item($r, 5).Value = $data[$i][3]
}

Write-Output "Updated $($data.Count) rows"